$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.207.50'
$ws.Range('E2').Value = '  -2.95%  '
$ws.Range('D3').Value = '1.959.25'
$ws.Range('E3').Value = '  -4.03%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.06'
$ws.Range('E5').Value = '  -3.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.618'
$ws.Range('E6').Value = '  -4.82%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '57.29'
$ws.Range('E7').Value = '  -12.80%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  -8.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '55.86'
$ws.Range('E10').Value = '  -5.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0852'
$ws.Range('E11').Value = '  +4.55%  '
$ws.Range('E12').Value = '  -0.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.12'
$ws.Range('E13').Value = '  -6.56%  '
$ws.Range('E14').Value = '  -9.32%  '
$ws.Range('D15').Value = '2.247.82'
$ws.Range('E15').Value = '  -4.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.39'
$ws.Range('E16').Value = '  -9.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.32'
$ws.Range('E17').Value = '  -6.30%  '
$ws.Range('D18').Value = '1.963.98'
$ws.Range('E18').Value = '  -3.82%  '
$ws.Range('D19').Value = '36.147.61'
$ws.Range('E19').Value = '  -2.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.76'
$ws.Range('E20').Value = '  -3.25%  '
$ws.Range('E21').Value = '  -2.84%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '230.52'
$ws.Range('E22').Value = '  -3.38%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.11'
$ws.Range('E23').Value = '  -6.56%  '
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.49'
$ws.Range('E25').Value = '  -4.33%  '
$ws.Range('E26').Value = '  -5.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.55'
$ws.Range('E27').Value = '  -4.98%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '166.41'
$ws.Range('E28').Value = '  +2.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.81'
$ws.Range('E29').Value = '  -1.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.127'
$ws.Range('E30').Value = '  -0.72%  '
$ws.Range('E31').Value = '  -3.72%  '
$ws.Range('E32').Value = '  -3.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.72'
$ws.Range('E33').Value = '  -8.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0639'
$ws.Range('E34').Value = '  +1.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.34'
$ws.Range('E35').Value = '  -6.84%  '
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.92'
$ws.Range('E38').Value = '  -7.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.17'
$ws.Range('E39').Value = '  -8.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.89'
$ws.Range('E40').Value = '  -4.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0958'
$ws.Range('E41').Value = '  -5.64%  '
$ws.Range('E42').Value = '  -5.28%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.18'
$ws.Range('E43').Value = '  -8.92%  '
$ws.Range('E44').Value = '  -4.27%  '
$ws.Range('E45').Value = '  -8.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '15.65'
$ws.Range('E46').Value = '  -9.73%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '88.38'
$ws.Range('E47').Value = '  -7.18%  '
$ws.Range('D48').Value = '1.342.45'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.25'
$ws.Range('E49').Value = '  -7.09%  '
$ws.Range('E50').Value = '  -4.49%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '44.48'
$ws.Range('E51').Value = '  -5.04%  '
